$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "22.418.70"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").Value = "1.572.51"
$ws.Range("E3").Value = "  +0.23%  "

# Row 4
$ws.Range("E4").Value = "  +0.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3765"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.18%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.35%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3423"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.46%  "

# Row 10
$ws.Range("E10").Value = "  -0.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07652"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.91%  "

# Row 12
$ws.Range("E12").Value = "  +0.23%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.998"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.80%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.924"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.75%  "

# Row 16
$ws.Range("D16").Value = "1.571.63"
$ws.Range("E16").Value = "  +0.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.63%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06759"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.70%  "

# Row 20
$ws.Range("E20").Value = "  +0.14%  "

# Row 21
$ws.Range("E21").Value = "  +2.32%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.210"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.52%  "

# Row 24
$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D24").Value = "22.411.64"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.423"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.75%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.723"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.97%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.12%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "146.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.032"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.35%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.87%  "

# Row 31
$ws.Range("D31").Value = "1.747.38"
$ws.Range("E31").Value = "  +0.08%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.189"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.36%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.009"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.02%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9947"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.33%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.46%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08563"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.62%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02543"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.15%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2313"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.55%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06582"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.89%  "

# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.337"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.80%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.454"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.36%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6444"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.95%  "

# Row 43
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.34%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.13%  "

# Row 45
$ws.Range("E45").Value = "  +0.23%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.797"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6011"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.314"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.61%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.083"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.94%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.99%  "

# Row 51
$ws.Range("E51").Value = "  +0.78%  "
